# Test_Cases.xlsx update: append 14 new test-case rows (78-91) to Sheet1,
# matching "Updated the spreadsheets with new tests".
#
# Columns: A=Test Case ID(opt) B=Test Case Name C=Status D=Key Requirement ID
#          E=Priority F=Spin Version G=Input Data Type H=JIRA Issue Number I=Negative Test?

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Step 1: write every brand-new text value (never used anywhere else on the
# sheet) in the exact sequence they need to be introduced, so that the
# resulting shared-string table is built up in the same order Excel produced.
# (Re-used/older values, e.g. "Passed"/"High"/"text"/"byte", are filled in
# afterwards in step 2 - doing so does not change their existing shared
# string id.)
# ---------------------------------------------------------------------------
$ws.Range("B78").Value = "property_scoping_02"
$ws.Range("D78").Value = "DFDL-8-018R"
$ws.Range("H78").Value = "DFDL-69/DFDL-70"

$ws.Range("D79").Value = "DFDL-7-028R"
$ws.Range("B79").Value = "defineFormat_01"

$ws.Range("B80").Value = "property_scoping_04"

$ws.Range("D81").Value = "DFDL-8-022R"
$ws.Range("D80").Value = "DFDL-8-021R"
$ws.Range("H80").Value = "DFDL-69/DFDL-71"
$ws.Range("B81").Value = "property_scoping_05"
$ws.Range("H81").Value = "DFDL-69/DFDL-72, DFDL-131"

$ws.Range("B82").Value = "escapeSchemeSimple"
$ws.Range("D82").Value = "DFDL-7-079R"
$ws.Range("H82").Value = "DFDL-269"

$ws.Range("B83").Value = "DelimProp_04"
$ws.Range("B84").Value = "DelimProp_02"
$ws.Range("B85").Value = "SeqGrp_02"
$ws.Range("B86").Value = "SeqGrp_03"
$ws.Range("B87").Value = "DelimProp_09"
$ws.Range("B88").Value = "DelimProp_10"
$ws.Range("B89").Value = "SeqGrp_04"
$ws.Range("B90").Value = "ParseSequence_4a"
$ws.Range("B91").Value = "AI000"
$ws.Range("D91").Value = "DFDL-12-038R"
$ws.Range("H91").Value = "DFDL-156"

# ---------------------------------------------------------------------------
# Step 2: fill in the remaining cells of each new row (these reuse strings
# already present in the workbook, plus the numeric "Spin Version" column).
# ---------------------------------------------------------------------------
$ws.Range("C78").Value = "Passed"
$ws.Range("E78").Value = "High"
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = "text"

$ws.Range("C79").Value = "Passed"
$ws.Range("E79").Value = "High"
$ws.Range("F79").Value = 3
$ws.Range("G79").Value = "text"
$ws.Range("H79").Value = "DFDL-69/DFDL-70"

$ws.Range("C80").Value = "Passed"
$ws.Range("E80").Value = "High"
$ws.Range("F80").Value = 3
$ws.Range("G80").Value = "byte"

$ws.Range("C81").Value = "Passed"
$ws.Range("E81").Value = "High"
$ws.Range("F81").Value = 3
$ws.Range("G81").Value = "byte"

$ws.Range("C82").Value = "Passed"
$ws.Range("E82").Value = "High"
$ws.Range("F82").Value = 3
$ws.Range("G82").Value = "text"

$ws.Range("C83").Value = "Passed"
$ws.Range("D83").Value = "DFDL-12-033R"
$ws.Range("E83").Value = "High"
$ws.Range("F83").Value = 3
$ws.Range("G83").Value = "text"
$ws.Range("H83").Value = "DFDL-269"

$ws.Range("C84").Value = "Passed"
$ws.Range("D84").Value = "DFDL-12-033R"
$ws.Range("E84").Value = "High"
$ws.Range("F84").Value = 3
$ws.Range("G84").Value = "text"
$ws.Range("H84").Value = "DFDL-269"

$ws.Range("C85").Value = "Passed"
$ws.Range("D85").Value = "DFDL-14-008R"
$ws.Range("E85").Value = "High"
$ws.Range("F85").Value = 3
$ws.Range("G85").Value = "text"
$ws.Range("H85").Value = "DFDL-269"
$ws.Range("I85").Value = "Yes"

$ws.Range("C86").Value = "Passed"
$ws.Range("D86").Value = "DFDL-14-008R"
$ws.Range("E86").Value = "High"
$ws.Range("F86").Value = 3
$ws.Range("G86").Value = "text"
$ws.Range("H86").Value = "DFDL-269"

$ws.Range("C87").Value = "Passed"
$ws.Range("D87").Value = "DFDL-12-033R"
$ws.Range("E87").Value = "High"
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = "text"
$ws.Range("H87").Value = "DFDL-269"

$ws.Range("C88").Value = "Passed"
$ws.Range("D88").Value = "DFDL-12-033R"
$ws.Range("E88").Value = "High"
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = "text"
$ws.Range("H88").Value = "DFDL-269"

$ws.Range("C89").Value = "Passed"
$ws.Range("D89").Value = "DFDL-14-008R"
$ws.Range("E89").Value = "High"
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = "text"
$ws.Range("H89").Value = "DFDL-269"

$ws.Range("C90").Value = "Passed"
$ws.Range("D90").Value = "DFDL-12-032R"
$ws.Range("E90").Value = "High"
$ws.Range("F90").Value = 3
$ws.Range("G90").Value = "text"
$ws.Range("H90").Value = "DFDL-269"

$ws.Range("C91").Value = "Passed"
$ws.Range("E91").Value = "High"
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = "text"

# ---------------------------------------------------------------------------
# Step 3: match the final selection / scroll position recorded in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("I86").Select()
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 2
